# Apply "contingencies with rene fine": add line7 and line8 rows, and update
# the C/D/E values of the lines/extr rows (the extr rows shift down by 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the new rows (8,9 and 16,17) inherit the same formatting
# (bold/bordered/centered column A) that the existing data rows already use,
# by copying the format from an existing fully-styled row before writing values.
$ws.Range("A7").Copy($ws.Range("A8:A9"))
$ws.Range("A15").Copy($ws.Range("A16:A17"))

# Final values for every data row (A..E), rows 2-17.
$data = @(
    @(2,  0, "line1", 7,  9,  $false),
    @(3,  1, "line2", 9,  8,  $true),
    @(4,  2, "line3", 8,  10, $true),
    @(5,  3, "line4", 8,  11, $true),
    @(6,  4, "line5", 10, 5,  $true),
    @(7,  5, "line6", 12, 8,  $true),
    @(8,  6, "line7", 14, 11, $true),
    @(9,  7, "line8", 16, 9,  $true),
    @(10, 8, "extr1", 5,  12, $true),
    @(11, 9, "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $false),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $false),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}
